$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated voltage magnitude (vm_pu) results for the "380 kV" case (rows 2-25)
$data = @{
    2 = @{ B = 1.02; C = 1.035734047817174; D = 1.043269584451086; E = 1.034750445191287; F = 1.051918106173402; I = 1.036959400417701; J = 1.040846206301867; K = 1.046043925429537; L = 1.037549022485983; M = 1.054668278148881; N = 1.042324327301544 }
    3 = @{ B = 1.02; C = 1.036703838607883; D = 1.044140052608532; E = 1.035574561430758; F = 1.052886528853457; I = 1.037130628279976; J = 1.041459566250191; K = 1.04672523875915; L = 1.038182351663187; M = 1.055449042855943; N = 1.042938558291349 }
    4 = @{ B = 1.02; C = 1.037331431404833; D = 1.044703656870942; E = 1.036108270677744; F = 1.053513530400893; I = 1.037239477210543; J = 1.04185591992719; K = 1.04716580836074; L = 1.038591961601045; M = 1.055953997324652; N = 1.043335474836031 }
    5 = @{ B = 1.02; C = 1.037595288231778; D = 1.044940679504489; E = 1.036332748822161; F = 1.053777208639293; I = 1.037284770540941; J = 1.042022418970643; K = 1.047350954490662; L = 1.038764113687836; M = 1.056166218690577; N = 1.043502210327223 }
    6 = @{ B = 1.02; C = 1.037639591918701; D = 1.044980481504938; E = 1.036370445917492; F = 1.053821486431562; I = 1.037292348100166; J = 1.042050367358927; K = 1.047382037247618; L = 1.038793015943322; M = 1.056201847955363; N = 1.043530198405424 }
    7 = @{ B = 1.02; C = 1.037334957007522; D = 1.044706823652708; E = 1.036111269747625; F = 1.053517053342855; I = 1.037240084257124; J = 1.041858145199775; K = 1.047168282566; L = 1.03859426209515; M = 1.055956833279876; N = 1.043337703268759 }
    8 = @{ B = 1.02; C = 1.036061777535873; D = 1.043563689429955; E = 1.035028865058131; F = 1.052245312408064; I = 1.037017670410177; J = 1.041053603642337; K = 1.046274237048945; L = 1.037763099161719; M = 1.054932192895021; N = 1.042532019170026 }
    9 = @{ B = 1.02; C = 1.033818866049183; D = 1.041552095777596; E = 1.033125034907463; F = 1.050007205839053; I = 1.036610866150621; J = 1.039631867915557; K = 1.04469666764552; L = 1.036297019801603; M = 1.053124757296032; N = 1.041108264415439 }
    10 = @{ B = 1.02; C = 1.032324031430063; D = 1.040212951513058; E = 1.03185823964231; F = 1.048517125783355; I = 1.036329694659536; J = 1.038681383555045; K = 1.043643564931695; L = 1.035318701800456; M = 1.051918590281251; N = 1.040156430258113 }
    11 = @{ B = 1.02; C = 1.031676863657248; D = 1.039633555608763; E = 1.03131029236648; F = 1.047872390158449; I = 1.036205587464062; J = 1.038269192235574; K = 1.043187241030139; L = 1.034894868963123; M = 1.051396031682126; N = 1.039743653579689 }
    12 = @{ B = 1.02; C = 1.031436493119729; D = 1.039418412665187; E = 1.03110684908114; F = 1.047632979586438; I = 1.036159134806601; J = 1.038115992820593; K = 1.043017694230859; L = 1.034737407129796; M = 1.051201888938192; N = 1.039590236603967 }
    13 = @{ B = 1.02; C = 1.031488052665707; D = 1.039464558350587; E = 1.031150484323998; F = 1.047684330649859; I = 1.036169115056973; J = 1.038148858821077; K = 1.043054064744164; L = 1.034771184637651; M = 1.05124353506636; N = 1.039623149277941 }
    14 = @{ B = 1.02; C = 1.03165699421619; D = 1.039615770382003; E = 1.031293473864775; F = 1.047852598912037; I = 1.03620175488547; J = 1.038256530621248; K = 1.043173227208914; L = 1.0348818537557; M = 1.051379984617891; N = 1.039730973984419 }
    15 = @{ B = 1.02; C = 1.031761086869003; D = 1.03970894647668; E = 1.031381586229814; F = 1.047956284237816; I = 1.036221818506388; J = 1.038322858442326; K = 1.043246640828687; L = 1.034950036509011; M = 1.051464050255724; N = 1.03979739599861 }
    16 = @{ B = 1.02; C = 1.032366984092349; D = 1.040251413940658; E = 1.031894617482622; F = 1.048559924938219; I = 1.036337881610838; J = 1.038708726224; K = 1.043673842900502; L = 1.035346825735866; M = 1.051953264945507; N = 1.040183811756795 }
    17 = @{ B = 1.02; C = 1.032747076036091; D = 1.040591813970045; E = 1.032216585396428; F = 1.048938701717477; I = 1.036410053998169; J = 1.038950603997931; K = 1.043941729446079; L = 1.035595664228877; M = 1.052260061983759; N = 1.040426033024908 }
    18 = @{ B = 1.02; C = 1.032968787377537; D = 1.040790408094622; E = 1.032404440054316; F = 1.049159681941829; I = 1.036451923372676; J = 1.039091626827566; K = 1.044097951870649; L = 1.035740786743745; M = 1.05243898439836; N = 1.040567256123136 }
    19 = @{ B = 1.02; C = 1.03304438691908; D = 1.040858131083645; E = 1.032468503141301; F = 1.049235038282871; I = 1.036466161132599; J = 1.039139701703383; K = 1.044151214352313; L = 1.035790266223735; M = 1.052499987692595; N = 1.040615399270791 }
    20 = @{ B = 1.02; C = 1.032706294714716; D = 1.040555287653292; E = 1.032182035463316; F = 1.048898057774253; I = 1.036402334117366; J = 1.038924659042244; K = 1.043912990982948; L = 1.035568968355809; M = 1.05222714836229; N = 1.040400051224408 }
    21 = @{ B = 1.02; C = 1.031607244683592; D = 1.039571240253475; E = 1.031251364550931; F = 1.047803046081123; I = 1.036192153029442; J = 1.038224826526763; K = 1.043138138145658; L = 1.034849265320397; M = 1.05133980475635; N = 1.039699224866482 }
    22 = @{ B = 1.02; C = 1.030916323393028; D = 1.038952938949236; E = 1.030666728704966; F = 1.047114991303483; I = 1.036057957599083; J = 1.037784275149836; K = 1.042650681347025; L = 1.034396577914845; M = 1.050781658608249; N = 1.039258047856069 }
    23 = @{ B = 1.02; C = 1.031282584648679; D = 1.039280673121998; E = 1.030976606159373; F = 1.047479701775225; I = 1.036129290911085; J = 1.038017870605332; K = 1.042909117316122; L = 1.03463657301316; M = 1.051077564702067; N = 1.039491975043906 }
    24 = @{ B = 1.02; C = 1.032724722007056; D = 1.040571792185989; E = 1.032197646917264; F = 1.048916422880177; I = 1.03640582310234; J = 1.038936382637528; K = 1.043925976754038; L = 1.035581031135946; M = 1.052242020694868; N = 1.040411791468542 }
    25 = @{ B = 1.02; C = 1.034398637900527; D = 1.042071807669346; E = 1.033616798209134; F = 1.05058546390287; I = 1.036717794577026; J = 1.039999893221264; K = 1.045104756021899; L = 1.036676204519527; M = 1.053592240329636; N = 1.041476812359293 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

Write-Output "Updated vm_pu values for rows 2-25 (380 kV case)"